# Auto-generated Excel COM-interop script
# Applies scheduled market-price/profit refresh values to the Leve profit sheets
# (columns H:N) across all 8 worksheets, per the authoritative diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 33
$ws.Range("H33").Value = 730.9091
$ws.Range("J33").Value = 350.25
$ws.Range("L33").Value = 350.25
$ws.Range("N33").Value = -808.25

# Row 62
$ws.Range("H62").Value = 1679.3
$ws.Range("I62").Value = 1421.4445
$ws.Range("K62").Value = 1421.4445
$ws.Range("M62").Value = -797.4445000000001

# Row 65
$ws.Range("H65").Value = 1679.3
$ws.Range("I65").Value = 1421.4445
$ws.Range("K65").Value = 7107.2225
$ws.Range("M65").Value = -3987.2225

# Row 92
$ws.Range("H92").Value = 1043.3572
$ws.Range("I92").Value = 394.44446
$ws.Range("K92").Value = 394.44446
$ws.Range("M92").Value = 853.5555400000001

# Row 106
$ws.Range("H106").Value = 2071
$ws.Range("I106").Value = 2062.9375
$ws.Range("K106").Value = 2062.9375
$ws.Range("M106").Value = -1431.9375

# Row 111
$ws.Range("H111").Value = 3378.875
$ws.Range("J111").Value = 3533.2856
$ws.Range("L111").Value = 10599.8568
$ws.Range("N111").Value = -16733.8568

# Row 112
$ws.Range("H112").Value = 3941.077
$ws.Range("J112").Value = 3936.1667
$ws.Range("L112").Value = 11808.5001
$ws.Range("N112").Value = -14024.5001

# Row 121
$ws.Range("H121").Value = 1297.2
$ws.Range("J121").Value = 1297.25
$ws.Range("L121").Value = 3891.75
$ws.Range("N121").Value = -7385.75

# Row 125
$ws.Range("H125").Value = 16662.375
$ws.Range("I125").Value = 26230
$ws.Range("J125").Value = 716.3333
$ws.Range("K125").Value = 236070
$ws.Range("L125").Value = 6446.9997
$ws.Range("M125").Value = -233610
$ws.Range("N125").Value = -11366.9997

# Row 132
$ws.Range("H132").Value = 1116.9231
$ws.Range("I132").Value = 960
$ws.Range("K132").Value = 2880
$ws.Range("M132").Value = -350

# Row 137
$ws.Range("H137").Value = 1555.2222
$ws.Range("I137").Value = 1366.2667
$ws.Range("K137").Value = 4098.800099999999
$ws.Range("M137").Value = -1548.800099999999


$ws = $wb.Worksheets.Item("ARM")

# Row 45
$ws.Range("H45").Value = 3268.7058
$ws.Range("I45").Value = 1892.4166
$ws.Range("J45").Value = 6571.8
$ws.Range("K45").Value = 1892.4166
$ws.Range("L45").Value = 6571.8
$ws.Range("M45").Value = -1515.4166
$ws.Range("N45").Value = -7325.8

# Row 102
$ws.Range("H102").Value = 5133.273
$ws.Range("I102").Value = 4056.0588
$ws.Range("K102").Value = 4056.0588
$ws.Range("M102").Value = -2434.0588

# Row 122
$ws.Range("H122").Value = 3205.4
$ws.Range("I122").Value = 2924.5833
$ws.Range("K122").Value = 8773.749899999999
$ws.Range("M122").Value = -6323.749899999999

# Row 128
$ws.Range("H128").Value = 43666.668
$ws.Range("J128").Value = 43666.668
$ws.Range("L128").Value = 43666.668
$ws.Range("N128").Value = -53626.668


$ws = $wb.Worksheets.Item("BSM")

# Row 134
$ws.Range("H134").Value = 3748.9119
$ws.Range("I134").Value = 3807.9697
$ws.Range("K134").Value = 11423.9091
$ws.Range("M134").Value = -8888.909100000001


$ws = $wb.Worksheets.Item("CRP")

# Row 6
$ws.Range("H6").Value = 5000
$ws.Range("I6").Value = 5000
$ws.Range("K6").Value = 5000
$ws.Range("M6").Value = -4887

# Row 31
$ws.Range("H31").Value = 4266.0356
$ws.Range("I31").Value = 3341.1765
$ws.Range("J31").Value = 5695.364
$ws.Range("K31").Value = 3341.1765
$ws.Range("L31").Value = 5695.364
$ws.Range("M31").Value = -3046.1765
$ws.Range("N31").Value = -6285.364

# Row 34
$ws.Range("H34").Value = 4266.0356
$ws.Range("I34").Value = 3341.1765
$ws.Range("J34").Value = 5695.364
$ws.Range("K34").Value = 3341.1765
$ws.Range("L34").Value = 5695.364
$ws.Range("M34").Value = -3139.1765
$ws.Range("N34").Value = -6099.364

# Row 134
$ws.Range("H134").Value = 4876.355
$ws.Range("I134").Value = 4042.7036
$ws.Range("K134").Value = 12128.1108
$ws.Range("M134").Value = -9593.110799999999

# Row 135
$ws.Range("H135").Value = 81590
$ws.Range("J135").Value = 81590
$ws.Range("L135").Value = 81590
$ws.Range("N135").Value = -91730


$ws = $wb.Worksheets.Item("CUL")

# Row 2
$ws.Range("H2").Value = 12309.818
$ws.Range("J2").Value = 16901
$ws.Range("L2").Value = 101406
$ws.Range("N2").Value = -101632

# Row 17
$ws.Range("H17").Value = 74
$ws.Range("J17").Value = 88.8
$ws.Range("L17").Value = 266.4
$ws.Range("N17").Value = -604.4

# Row 26
$ws.Range("H26").Value = 285943.44
$ws.Range("I26").Value = 500076.25
$ws.Range("J26").Value = 433
$ws.Range("K26").Value = 1500228.75
$ws.Range("L26").Value = 1299
$ws.Range("M26").Value = -1499940.75
$ws.Range("N26").Value = -1875

# Row 33
$ws.Range("H33").Value = 505.55554
$ws.Range("I33").Value = 243
$ws.Range("K33").Value = 1458
$ws.Range("M33").Value = -1175

# Row 69
$ws.Range("H69").Value = 646.25
$ws.Range("J69").Value = 642.5
$ws.Range("L69").Value = 1927.5
$ws.Range("N69").Value = -3549.5

# Row 72
$ws.Range("H72").Value = 646.25
$ws.Range("J72").Value = 642.5
$ws.Range("L72").Value = 5782.5
$ws.Range("N72").Value = -13894.5

# Row 80
$ws.Range("H80").Value = 3247.25
$ws.Range("J80").Value = 3247.25
$ws.Range("L80").Value = 9741.75
$ws.Range("N80").Value = -11613.75

# Row 82
$ws.Range("H82").Value = 5001265.5
$ws.Range("I82").Value = 5001265.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 15003796.5
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -15003390.5

# Row 83
$ws.Range("H83").Value = 3247.25
$ws.Range("J83").Value = 3247.25
$ws.Range("L83").Value = 29225.25
$ws.Range("N83").Value = -38585.25

# Row 85
$ws.Range("H85").Value = 5001265.5
$ws.Range("I85").Value = 5001265.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 15003796.5
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -15002392.5

# Row 119
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()


$ws = $wb.Worksheets.Item("GSM")

# Row 97
$ws.Range("H97").Value = 1180.5834
$ws.Range("I97").Value = 849.7143
$ws.Range("J97").Value = 1643.8
$ws.Range("K97").Value = 849.7143
$ws.Range("L97").Value = 1643.8
$ws.Range("M97").Value = -353.7143
$ws.Range("N97").Value = -2635.8

# Row 113
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830

# Row 136
$ws.Range("H136").Value = 37905.75
$ws.Range("J136").Value = 37905.75
$ws.Range("L136").Value = 113717.25
$ws.Range("N136").Value = -118817.25


$ws = $wb.Worksheets.Item("LTW")

# Row 68
$ws.Range("H68").Value = 1991.75
$ws.Range("J68").Value = 1989.3334
$ws.Range("L68").Value = 1989.3334
$ws.Range("N68").Value = -3487.3334

# Row 71
$ws.Range("H71").Value = 1991.75
$ws.Range("J71").Value = 1989.3334
$ws.Range("L71").Value = 9946.666999999999
$ws.Range("N71").Value = -17434.667

# Row 132
$ws.Range("H132").Value = 9818.303
$ws.Range("J132").Value = 5705.091
$ws.Range("L132").Value = 17115.273
$ws.Range("N132").Value = -22175.273

# Row 136
$ws.Range("H136").Value = 4882.087
$ws.Range("J136").Value = 5878.6
$ws.Range("L136").Value = 17635.8
$ws.Range("N136").Value = -22735.8


$ws = $wb.Worksheets.Item("WVR")

# Row 5
$ws.Range("H5").Value = 11996.3
$ws.Range("I5").Value = 10001
$ws.Range("J5").Value = 12218
$ws.Range("K5").Value = 10001
$ws.Range("L5").Value = 12218
$ws.Range("M5").Value = -9889
$ws.Range("N5").Value = -12442

# Row 14
$ws.Range("H14").Value = 16773.258
$ws.Range("I14").Value = 14410.883
$ws.Range("K14").Value = 14410.883
$ws.Range("M14").Value = -14242.883

# Row 122
$ws.Range("H122").Value = 4071.04
$ws.Range("I122").Value = 3898.3333
$ws.Range("J122").Value = 4515.143
$ws.Range("K122").Value = 11694.9999
$ws.Range("L122").Value = 13545.429
$ws.Range("M122").Value = -9244.999899999999
$ws.Range("N122").Value = -18445.429

# Row 136
$ws.Range("H136").Value = 2661.7026
$ws.Range("I136").Value = 2024.9697
$ws.Range("K136").Value = 6074.909100000001
$ws.Range("M136").Value = -3524.909100000001

